# Reformat the date/time columns that came back from the new SQLAlchemy DB
# read: the raw "YYYY-MM-DD-HHMMSS" strings are replaced with a human
# readable "DD MON YYYY HH:MM" style layout, and the affected columns are
# forced to text formatting so Excel never reinterprets them as real dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A ("dateTime") blocks --------------------------------------
# Each block shares one dateTime value across its rows.
$blocksA = @(
    @{ Rows = 2..5;   Value = "11 OCT 2019  15:13" },
    @{ Rows = 6..11;  Value = "11 OCT 2019  17:33" },
    @{ Rows = 12..17; Value = "11 OCT 2019  16:23" },
    @{ Rows = 18..23; Value = "09 OCT 2019  10:00" }
)
foreach ($block in $blocksA) {
    foreach ($r in $block.Rows) {
        $ws.Cells.Item($r, 1).Value = $block.Value
    }
}

# --- Column M ("onset_date_time") blocks --------------------------------
$blocksM = @(
    @{ Rows = 2..5;   Value = "11 OCT 2019-15:12" },
    @{ Rows = 6..11;  Value = "11 OCT 2019-17:32" },
    @{ Rows = 12..17; Value = "11 OCT 2019-16:22" },
    @{ Rows = 18..23; Value = "11 OCT 2019-17:32" }
)
foreach ($block in $blocksM) {
    foreach ($r in $block.Rows) {
        $ws.Cells.Item($r, 13).Value = $block.Value
    }
}

# --- Column S ("tsunami_arival_date"), only populated for the Tsunami rows
foreach ($r in 2..5) {
    $ws.Cells.Item($r, 19).Value = "11 OCT 2019 17:32"
}

# --- Column V ("typhoon_arival_date"), only populated for the Typhoon rows
foreach ($r in 12..17) {
    $ws.Cells.Item($r, 22).Value = "11 OCT 2019 17:32"
}

# --- Force text formatting on the reformatted date/time columns so Excel
#     keeps storing them verbatim instead of coercing to a date serial.
#     Only touch cells that actually hold data - touching the whole column
#     would materialise blank styled cells that should stay absent.
$ws.Range("A2:A23").NumberFormat = "@"
$ws.Range("M1:M23").NumberFormat = "@"
$ws.Range("S1:S5").NumberFormat = "@"
$ws.Range("V1").NumberFormat = "@"
$ws.Range("V12:V17").NumberFormat = "@"

# --- Restore the active selection left by the author after the edit.
$ws.Range("A23").Select() | Out-Null
